$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row - new outlier-detection columns
$ws.Range("F1").Value = "KNN_Outliers_MAD"
$ws.Range("G1").Value = "SVM_Outliers_MAD"
$ws.Range("H1").Value = "RF_Outliers_MAD"

# Copy the header style (bold, centered, bordered) from an existing header cell
$ws.Range("E1").Copy()
$ws.Range("F1:H1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Boolean outlier flags for each data row
$ws.Range("F2").Value = $false
$ws.Range("G2").Value = $false
$ws.Range("H2").Value = $false

$ws.Range("F3").Value = $false
$ws.Range("G3").Value = $false
$ws.Range("H3").Value = $false

$ws.Range("F4").Value = $false
$ws.Range("G4").Value = $false
$ws.Range("H4").Value = $false

$ws.Range("F5").Value = $true
$ws.Range("G5").Value = $false
$ws.Range("H5").Value = $false

$ws.Range("F6").Value = $false
$ws.Range("G6").Value = $false
$ws.Range("H6").Value = $false

$ws.Range("F7").Value = $false
$ws.Range("G7").Value = $false
$ws.Range("H7").Value = $false

$ws.Range("F8").Value = $false
$ws.Range("G8").Value = $false
$ws.Range("H8").Value = $false

$ws.Range("F9").Value = $false
$ws.Range("G9").Value = $false
$ws.Range("H9").Value = $false
